$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for two new columns in the "left" (data-entry) table:
#    firstCond / secondCond are inserted right after subNr (column A),
#    pushing the former B:G (firstSkill..) over to D:I.
# ---------------------------------------------------------------------------
$ws.Range("B1:C1").EntireColumn.Insert()

# 2. Make room for two new columns in the "right" (lookup/formula) table:
#    the firstCond/secondCond formulas are inserted right after the
#    subNr-concatenation column (now J), pushing the old I:L formulas
#    (now at J+1..) over to M:P.
$ws.Range("K1:L1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 3. Header row. NB: shared strings get interned in the order they are
#    first written, so these particular cells are written in this
#    particular order to match the original author's string table.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "firstCond"
$ws.Range("F1").Value = "firstScore"
$ws.Range("G1").Value = "secondScore"

# ---------------------------------------------------------------------------
# 4. Condition values for each participant (alternating Teaching/Performing)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Teaching"
$ws.Range("C2").Value = "Performing"

$ws.Range("C1").Value = "secondCond"

$ws.Range("B3").Value = "Performing"
$ws.Range("C3").Value = "Teaching"

$ws.Range("B4").Value = "Teaching"
$ws.Range("C4").Value = "Performing"

$ws.Range("B5").Value = "Performing"
$ws.Range("C5").Value = "Teaching"

$ws.Range("B6").Value = "Teaching"
$ws.Range("C6").Value = "Performing"

$ws.Range("B7").Value = "Performing"
$ws.Range("C7").Value = "Teaching"

$ws.Range("B8").Value = "Teaching"
$ws.Range("C8").Value = "Performing"

$ws.Range("B9").Value = "Performing"
$ws.Range("C9").Value = "Teaching"

# ---------------------------------------------------------------------------
# 5. Mirror the new conditions into the lookup table's header + formulas
#    (K = firstCond, L = secondCond), matching the existing firstSkill/
#    secondSkill lookup-formula pattern already used in that table.
# ---------------------------------------------------------------------------
$ws.Range("K1").Value = "firstCond"
$ws.Range("L1").Value = "secondCond"

# Re-assert the pre-existing lookup formulas (M:P) before the new K/L ones
# so everything regroups into the expected shared-formula ranges.
$ws.Range("M3:M9").Formula = "=D3"
$ws.Range("N3:N9").Formula = "=E3"
$ws.Range("O3:O9").Formula = "=F3"
$ws.Range("P3:P9").Formula = "=CONCATENATE(G3,I3)"

$ws.Range("K2").Formula = "=B2"
$ws.Range("L2").Formula = "=C2"
$ws.Range("K3:K9").Formula = "=B3"
$ws.Range("L3:L9").Formula = "=C3"

# ---------------------------------------------------------------------------
# 6. Selection left however the author last left it.
# ---------------------------------------------------------------------------
$ws.Range("J15").Select()
